{"js": "// The paragraph \"<id>p136v_1</id>\" is currently split across three runs:\n//   1) \"<id>\"   (Courier New, color 7f6000, size 9pt/sz 18)\n//   2) \"p136v_1\" (plain, color black)\n//   3) \"</id>\"  (Courier New, color 7f6000, size 9pt/sz 18)\n// The edit merges them into a single run containing the full text\n// \"<id>p136v_1</id>\", taking on the formatting of the first run.\nconst body = context.document.body;\n\nconst results = body.search(\"<id>p136v_1</id>\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text '<id>p136v_1</id>' not found in document body.\");\n}\n\n// Replacing the whole matched range with its own text collapses the three\n// runs into a single run, inheriting the formatting of the first run in the\n// range (matches the formatting kept in the target OOXML: Courier New,\n// color 7f6000, sz 18).\nconst target = results.items[0];\ntarget.insertText(\"<id>p136v_1</id>\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The paragraph \"<id>p136v_1</id>\" is currently split across three runs:\n#   1) \"<id>\"    (Courier New, color 7f6000, sz 18)\n#   2) \"p136v_1\" (plain run, color 000000)\n#   3) \"</id>\"   (Courier New, color 7f6000, sz 18)\n# The edit merges them into a single run containing the full text\n# \"<id>p136v_1</id>\", keeping the formatting of the first run (\"<id>\").\n$d = $word.ActiveDocument\n\n# Locate the whole combined text and the first fragment \"<id>\" so we know\n# exactly where run 1 ends - that's the split point between run 1 and the\n# content (runs 2 & 3) that needs to be folded into it.\n$full = $d.Content\n$full.Find.Execute(\"<id>p136v_1</id>\") | Out-Null\n$fullStart = $full.Start\n$fullEnd = $full.End\n\n$firstRun = $d.Content\n$firstRun.Find.Execute(\"<id>\") | Out-Null\n$splitPos = $firstRun.End\n\nif ($firstRun.Start -ne $fullStart) {\n    throw \"Unexpected range alignment while merging <id> runs.\"\n}\n\n$tailText = \"p136v_1</id>\"\n\n# Delete the trailing runs' text (runs 2 and 3) ...\n$tailRange = $d.Range($splitPos, $fullEnd)\n$tailRange.Delete()\n\n# ... then append it back onto run 1's own range so Word folds it into that\n# run and keeps run 1's character formatting (Courier New / 7f6000 / sz 18).\n$runOneRange = $d.Range($fullStart, $splitPos)\n$runOneRange.InsertAfter($tailText)\n"}
